# Insert a new row at position 151 (shifts existing rows 151..181 down to 152..182)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with its data
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 44543
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = 100112043
$ws.Range("G151").Value = "Pepino ensalada"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 300
$ws.Range("K151").Value = 7000
$ws.Range("L151").Value = 7500
$ws.Range("M151").Value = 7250
$ws.Range("N151").Value = "$/caja 80 unidades"
$ws.Range("O151").Value = "Región del Maule"
$ws.Range("P151").Value = 91
$ws.Range("Q151").Value = 80
$ws.Range("R151").Value = "Hortaliza"
